# Edit script: updates Betfair Back/Lay odds workbook for 2026-01-01
# - Inserts a new fixture row (Welsh Premiership: Llanelli Town v Penybont FC)
#   between the existing "Saudi 1st Division 09:35" row and "Saudi 1st Division 12:00" row.
# - Updates a number of odds values across several other fixture rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row at position 6 (pushes the old rows 6-11 down to 7-12).
#    Copy row 5 first so the new row inherits the same (default/General)
#    formatting as the other data rows, instead of Excel guessing formats.
# ---------------------------------------------------------------------------
$ws.Rows("5:5").Copy()
$ws.Rows("6:6").Insert()
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Populate the newly inserted row 6 with the new fixture's data.
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = "Welsh Premiership"

# Column B holds a plain text date string (e.g. "2026-01-01"). Setting it
# directly would make Excel auto-convert it to a real date, so force a Text
# format first, assign the value, then clear the format again so the cell
# ends up with default (General) styling but a literal text value, exactly
# like the surrounding cells.
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "2026-01-01"
$ws.Range("B6").ClearFormats()

$ws.Range("C6").Value = "11:00:00"
$ws.Range("D6").Value = "Llanelli Town"
$ws.Range("E6").Value = "Penybont FC"

$ws.Range("F6").Value = 1.01
$ws.Range("G6").Value = 1000
$ws.Range("H6").Value = 1.01
$ws.Range("I6").Value = 1000
$ws.Range("J6").Value = 5.3
$ws.Range("K6").Value = 1000
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 0
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = 0
$ws.Range("T6").Value = 0
$ws.Range("U6").Value = 0
$ws.Range("V6").Value = 0
$ws.Range("W6").Value = 0
$ws.Range("X6").Value = 0
$ws.Range("Y6").Value = 0
$ws.Range("Z6").Value = 0
$ws.Range("AA6").Value = 0
$ws.Range("AB6").Value = 0
$ws.Range("AC6").Value = 0
$ws.Range("AD6").Value = 0
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 0
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 0
$ws.Range("AK6").Value = 0
$ws.Range("AL6").Value = 0
$ws.Range("AM6").Value = 0
$ws.Range("AN6").Value = 0
$ws.Range("AO6").Value = 0

# ---------------------------------------------------------------------------
# 3. Update changed odds values on other rows (positions after the insert).
# ---------------------------------------------------------------------------

# Row 2 - Australian A-League Men: Auckland FC v Newcastle Jets
$ws.Range("AC2").Value = 13
$ws.Range("AJ2").Value = 20
$ws.Range("AK2").Value = 19
$ws.Range("AN2").Value = 6.6

# Row 3 - Australian A-League Men: Western Sydney Wanderers v Macarthur FC
$ws.Range("AB3").Value = 12.5

# Row 4 - Welsh Premiership: Colwyn Bay v Flint Town United
$ws.Range("F4").Value = 1.79
$ws.Range("H4").Value = 3.75
$ws.Range("I4").Value = 5.3
$ws.Range("J4").Value = 3.35
$ws.Range("K4").Value = 4.4
$ws.Range("L4").Value = 1.29
$ws.Range("M4").Value = 1.05
$ws.Range("N4").Value = 4.1
$ws.Range("O4").Value = 1.25
$ws.Range("P4").Value = 2.06
$ws.Range("Q4").Value = 1.72
$ws.Range("R4").Value = 1.42
$ws.Range("S4").Value = 2.66
$ws.Range("T4").Value = 1.7
$ws.Range("U4").Value = 2.12
$ws.Range("V4").Value = 1.24
$ws.Range("W4").Value = 1.93
$ws.Range("Y4").Value = 24
$ws.Range("AB4").Value = 12.5
$ws.Range("AC4").Value = 11
$ws.Range("AG4").Value = 12.5

# Row 7 (was row 6) - Saudi 1st Division: Al Faisaly ( KSA ) v Al-Raed (KSA)
$ws.Range("F7").Value = 1.02
$ws.Range("H7").Value = 1.02
$ws.Range("J7").Value = 1.02
$ws.Range("N7").Value = 1.3
$ws.Range("P7").Value = 1.3
$ws.Range("Q7").Value = 1.02
$ws.Range("R7").Value = 1.18
$ws.Range("S7").Value = 1.36

# Row 8 (was row 7) - English Premier League: Crystal Palace v Fulham
$ws.Range("N8").Value = 3.55
$ws.Range("T8").Value = 1.89
$ws.Range("Y8").Value = 12.5
$ws.Range("AC8").Value = 7.2
$ws.Range("AL8").Value = 42
$ws.Range("AO8").Value = 55

# Row 9 (was row 8) - English Premier League: Liverpool v Leeds
$ws.Range("M9").Value = 1.05
$ws.Range("N9").Value = 5.1
$ws.Range("R9").Value = 1.54
$ws.Range("X9").Value = 22
$ws.Range("AA9").Value = 170
$ws.Range("AB9").Value = 10.5
$ws.Range("AD9").Value = 23
$ws.Range("AF9").Value = 10.5
$ws.Range("AG9").Value = 9.2
$ws.Range("AH9").Value = 19.5
$ws.Range("AK9").Value = 14.5
$ws.Range("AL9").Value = 28
$ws.Range("AO9").Value = 70

# Row 10 (was row 9) - Israeli Premier League: Beitar Jerusalem v Hapoel Tel Aviv
$ws.Range("I10").Value = 4.1

# Row 11 (was row 10) - English Premier League: Brentford v Tottenham
$ws.Range("F11").Value = 2.3
$ws.Range("G11").Value = 2.32
$ws.Range("I11").Value = 3.4
$ws.Range("O11").Value = 1.33
$ws.Range("P11").Value = 1.92
$ws.Range("Q11").Value = 2.04
$ws.Range("R11").Value = 1.36
$ws.Range("S11").Value = 3.6
$ws.Range("T11").Value = 1.82
$ws.Range("U11").Value = 2.18
$ws.Range("X11").Value = 14.5
$ws.Range("Y11").Value = 13
$ws.Range("AB11").Value = 10
$ws.Range("AD11").Value = 15
$ws.Range("AE11").Value = 40
$ws.Range("AF11").Value = 14.5
$ws.Range("AI11").Value = 55
$ws.Range("AL11").Value = 38
$ws.Range("AM11").Value = 85
$ws.Range("AN11").Value = 19.5
$ws.Range("AO11").Value = 38

# Row 12 (was row 11) - English Premier League: Sunderland v Man City
# No value changes for this row; it keeps all of its original data.

Write-Output "Edit complete"
